$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "txn_1753294317066_hbvn3z8to"
$ws.Range("B5").Value = "settlement_1753294282531_mmatj9jez"
$ws.Range("C5").Value = "user_1753125931723_8ftkkx2pf"
$ws.Range("D5").Value = 5
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "111"
$ws.Range("F5").Value = "payment-proof-1753294317020-14169049.webp"
$ws.Range("G5").Value = "completed"
$ws.Range("H5").Value = "INR"
$ws.Range("I5").Value = "2025-07-23T18:11:57.066Z"
